$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert a new first column (everything shifts from A..G to B..H)
# ---------------------------------------------------------------------------
$ws1.Columns.Item(1).Insert()

# ---------------------------------------------------------------------------
# 2. Insert the four "screen" marker rows (top-down, using the final target
#    row numbers since each insert pushes everything below it down by one)
# ---------------------------------------------------------------------------
$ws1.Rows.Item(11).Insert()
$ws1.Rows.Item(14).Insert()
$ws1.Rows.Item(17).Insert()
$ws1.Rows.Item(20).Insert()

# ---------------------------------------------------------------------------
# 3. Fix up the header row labels (display.text/display.hint renamed)
# ---------------------------------------------------------------------------
$ws1.Range("D1").Value2 = "display.prompt.text"
$ws1.Range("E1").Value2 = "display.hint.text"

# Give the new "clause" header cell (A1) the same plain style used by H1
# (font size 18, no fill/border) before overwriting its text.
$ws1.Range("H1").Copy($ws1.Range("A1"))
$ws1.Range("A1").Value2 = "clause"

# ---------------------------------------------------------------------------
# 4. Populate the marker rows in column A, matching the style of H1/A1 and
#    leaving the rest of the row completely empty
# ---------------------------------------------------------------------------
function Set-ScreenMarker($rowNum, $text) {
    $ws1.Range("H1").Copy($ws1.Range("A$rowNum"))
    $ws1.Range("A$rowNum").Value2 = $text
    $ws1.Range("B${rowNum}:H$rowNum").Style = "Normal"
    $ws1.Range("B${rowNum}:H$rowNum").ClearContents()
}

Set-ScreenMarker 11 "begin screen"
Set-ScreenMarker 14 "end screen"
Set-ScreenMarker 17 "begin screen"
Set-ScreenMarker 20 "end screen"

# ---------------------------------------------------------------------------
# 5. Column widths / styles: give new column A the plain default width/style
# ---------------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 10.83203125

# ---------------------------------------------------------------------------
# 6. Restore the selection to match the authored file
# ---------------------------------------------------------------------------
$ws1.Range("D18").Select()
